# Add a new "Sheet3" with a PGSQL-vs-SQLite timing comparison table, plus a
# clustered-column chart plotting it, matching the upstream commit.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet, named "Sheet3", positioned after Sheet2 ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet2)
$ws3.Name = "Sheet3"

# --- 2. Header row ---
$ws3.Range("A1").Value = "Quesies"
$ws3.Range("B1").Value = "PGSQL"
$ws3.Range("C1").Value = "SQLite"

# --- 3. Data rows: query number, PGSQL real time, SQLite real time ---
$data = @(
    @(1, 0.21, 37.551),
    @(2, 41.816, 1.487),
    @(3, 42.969, 11.024),
    @(4, 23.594, 1.213),
    @(5, 22.624, 6.351),
    @(6, 17.366, 5.234),
    @(7, 15.578, 22.046),
    @(8, 5.188, 35.84),
    @(9, 27.861, 94.626),
    @(10, 10.791, 5.337),
    @(11, 4.346, 3.109),
    @(12, 2.811, 5.461),
    @(13, 5.651, 78.3999),
    @(14, 27.018, 7.241),
    @(15, 8.102, 5.486),
    @(16, 13.233, 2.008),
    @(18, 14.61, 6.101),
    @(19, 2.249, 7.265),
    @(21, 4.761, 30.657)
)

$row = 2
foreach ($d in $data) {
    $ws3.Cells.Item($row, 1).Value = $d[0]
    $ws3.Cells.Item($row, 2).Value = $d[1]
    $ws3.Cells.Item($row, 3).Value = $d[2]
    $row = $row + 1
}

# --- 4. Chart: clustered column chart of PGSQL vs SQLite timings ---
$chartObj = $ws3.Shapes.AddChart2(201, 51).Chart
$chartObj.SetSourceData($ws3.Range("B1:C20"))
$chartObj.ChartType = 51
$chartObj.SeriesCollection(1).XValues = $ws3.Range("A2:A20")
$chartObj.SeriesCollection(2).XValues = $ws3.Range("A2:A20")
$chartObj.HasLegend = $true
$chartObj.Legend.Position = -4107
$chartObj.ChartGroups(1).GapWidth = 219
$chartObj.ChartGroups(1).Overlap = -27

# --- 5. Selection / view bookkeeping to mirror the saved workbook state ---
$sheet1.Activate()
$sheet1.Range("A33").Select()

$sheet2.Activate()
$sheet2.Range("B2:B22").Select()

$ws3.Activate()
$ws3.Range("D26").Select()
